# Append a new log entry (row 10) to the git commits log worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "2025-08-07 23:36"
$ws.Range("B10").Value = "0969bf5"
$ws.Range("C10").Value = "[UI Enhancement]: Major mobile navigation and typography improvements"
$ws.Range("D10").Value = 8
$ws.Range("E10").Value = 585
$ws.Range("F10").Value = 36
$ws.Range("G10").Value = "Major mobile UI enhancements: Heebo font implementation, mobile navigation redesign, header/sidebar layout optimization, dashboard spacing improvements"
$ws.Range("H10").Value = "Local"
